$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4742.3105
$ws.Range("I40").Value = 3811.75
$ws.Range("J40").Value = 4891.2
$ws.Range("K40").Value = 3811.75
$ws.Range("L40").Value = 4891.2
$ws.Range("M40").Value = -3636.75
$ws.Range("N40").Value = -5241.2

$ws.Range("H43").Value = 100007736
$ws.Range("I43").Value = 250000590
$ws.Range("J43").Value = 12499
$ws.Range("K43").Value = 250000590
$ws.Range("L43").Value = 12499
$ws.Range("M43").Value = -250000521
$ws.Range("N43").Value = -12637

$ws.Range("H62").Value = 6010.5557
$ws.Range("J62").Value = 8000
$ws.Range("L62").Value = 8000
$ws.Range("N62").Value = -9248

$ws.Range("H65").Value = 6010.5557
$ws.Range("J65").Value = 8000
$ws.Range("L65").Value = 40000
$ws.Range("N65").Value = -46240

$ws.Range("H96").Value = 319.42856
$ws.Range("I96").Value = 361.4
$ws.Range("K96").Value = 1084.2
$ws.Range("M96").Value = 288.8000000000002

$ws.Range("H116").Value = 6707.5557
$ws.Range("I116").Value = 4580.625
$ws.Range("K116").Value = 4580.625
$ws.Range("M116").Value = -1138.625

$ws.Range("H132").Value = 20002964
$ws.Range("I132").Value = 23258772
$ws.Range("J132").Value = 2998.4285
$ws.Range("K132").Value = 69776316
$ws.Range("L132").Value = 8995.2855
$ws.Range("M132").Value = -69773786
$ws.Range("N132").Value = -14055.2855

$ws.Range("H135").Value = 926.5333000000001
$ws.Range("I135").Value = 626.96
$ws.Range("K135").Value = 5642.64
$ws.Range("M135").Value = -3107.64

$ws.Range("H137").Value = 70740.38
$ws.Range("I137").Value = 163433.64
$ws.Range("J137").Value = 2765.3333
$ws.Range("K137").Value = 490300.92
$ws.Range("L137").Value = 8295.999899999999
$ws.Range("M137").Value = -487750.92
$ws.Range("N137").Value = -13395.9999

$ws.Range("H138").Value = 3144.9556
$ws.Range("I138").Value = 2002.7333
$ws.Range("J138").Value = 3716.0667
$ws.Range("K138").Value = 6008.199900000001
$ws.Range("L138").Value = 11148.2001
$ws.Range("M138").Value = -868.1999000000005
$ws.Range("N138").Value = -21428.2001

$ws.Range("H141").Value = 11962.206
$ws.Range("I141").Value = 6363.8335
$ws.Range("K141").Value = 19091.5005
$ws.Range("M141").Value = -13911.5005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2778931.2
$ws.Range("I2").Value = 4274122
$ws.Range("K2").Value = 4274122
$ws.Range("M2").Value = -4274009

$ws.Range("H32").Value = 6430.398
$ws.Range("I32").Value = 3396.9875
$ws.Range("J32").Value = 19912.223
$ws.Range("K32").Value = 3396.9875
$ws.Range("L32").Value = 19912.223
$ws.Range("M32").Value = -3109.9875
$ws.Range("N32").Value = -20486.223

$ws.Range("H45").Value = 5758196
$ws.Range("I45").Value = 11067734
$ws.Range("J45").Value = 6196.25
$ws.Range("K45").Value = 11067734
$ws.Range("L45").Value = 6196.25
$ws.Range("M45").Value = -11067357
$ws.Range("N45").Value = -6950.25

$ws.Range("H74").Value = 29170.973
$ws.Range("I74").Value = 4322.9355
$ws.Range("K74").Value = 4322.9355
$ws.Range("M74").Value = -3448.9355

$ws.Range("H77").Value = 29170.973
$ws.Range("I77").Value = 4322.9355
$ws.Range("K77").Value = 21614.6775
$ws.Range("M77").Value = -17246.6775

$ws.Range("H110").Value = 1158679.4
$ws.Range("I110").Value = 1208926.2
$ws.Range("K110").Value = 1208926.2
$ws.Range("M110").Value = -1206881.2

$ws.Range("H116").Value = 2778931.2
$ws.Range("I116").Value = 4274122
$ws.Range("K116").Value = 4274122
$ws.Range("M116").Value = -4271828

$ws.Range("H132").Value = 3443.4546
$ws.Range("I132").Value = 2222.3076
$ws.Range("J132").Value = 5207.3335
$ws.Range("K132").Value = 6666.9228
$ws.Range("L132").Value = 15622.0005
$ws.Range("M132").Value = -4136.9228
$ws.Range("N132").Value = -20682.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2778931.2
$ws.Range("I3").Value = 4274122
$ws.Range("K3").Value = 4274122
$ws.Range("M3").Value = -4274008

$ws.Range("H63").Value = 24500
$ws.Range("I63").Value = 15000
$ws.Range("K63").Value = 15000
$ws.Range("M63").Value = -14314

$ws.Range("H66").Value = 24500
$ws.Range("I66").Value = 15000
$ws.Range("K66").Value = 45000
$ws.Range("M66").Value = -41568

$ws.Range("H107").Value = 5496109
$ws.Range("I107").Value = 5953826.5
$ws.Range("K107").Value = 5953826.5
$ws.Range("M107").Value = -5951906.5

$ws.Range("H134").Value = 4254.4136
$ws.Range("I134").Value = 1414.7222
$ws.Range("K134").Value = 4244.1666
$ws.Range("M134").Value = -1709.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1492.8462
$ws.Range("I16").Value = 1051.5
$ws.Range("K16").Value = 1051.5
$ws.Range("M16").Value = -764.5

$ws.Range("H22").Value = 452.26086
$ws.Range("I22").Value = 448.41666
$ws.Range("K22").Value = 448.41666
$ws.Range("M22").Value = -98.41665999999998

$ws.Range("H105").Value = 1166.7
$ws.Range("I105").Value = 723.8570999999999
$ws.Range("K105").Value = 723.8570999999999
$ws.Range("M105").Value = 1023.1429

$ws.Range("H107").Value = 1454.4595
$ws.Range("J107").Value = 1153.7142
$ws.Range("L107").Value = 1153.7142
$ws.Range("N107").Value = -4993.7142

$ws.Range("H113").Value = 1492.8462
$ws.Range("I113").Value = 1051.5
$ws.Range("K113").Value = 1051.5
$ws.Range("M113").Value = 1118.5

$ws.Range("H132").Value = 52067.684
$ws.Range("I132").Value = 36202.242
$ws.Range("J132").Value = 103189.664
$ws.Range("K132").Value = 108606.726
$ws.Range("L132").Value = 309568.992
$ws.Range("M132").Value = -106076.726
$ws.Range("N132").Value = -314628.992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 5424.8887
$ws.Range("I14").Value = 5424.8887
$ws.Range("K14").Value = 16274.6661
$ws.Range("M14").Value = -16101.6661

$ws.Range("H98").Value = 1558.2727
$ws.Range("J98").Value = 1732.3334
$ws.Range("L98").Value = 5197.0002
$ws.Range("N98").Value = -8193.0002

$ws.Range("H131").Value = 16031808
$ws.Range("I131").Value = 10419416
$ws.Range("J131").Value = 18526204
$ws.Range("K131").Value = 31258248
$ws.Range("L131").Value = 55578612
$ws.Range("M131").Value = -31253208
$ws.Range("N131").Value = -55588692

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 628.5
$ws.Range("I107").Value = 560
$ws.Range("K107").Value = 560
$ws.Range("M107").Value = 1360

$ws.Range("H122").Value = 264560.53
$ws.Range("I122").Value = 319825.47
$ws.Range("J122").Value = 6657.5
$ws.Range("K122").Value = 959476.4099999999
$ws.Range("L122").Value = 19972.5
$ws.Range("M122").Value = -957026.4099999999
$ws.Range("N122").Value = -24872.5

$ws.Range("H132").Value = 3221.639
$ws.Range("I132").Value = 3157.8667
$ws.Range("J132").Value = 3540.5
$ws.Range("K132").Value = 9473.6001
$ws.Range("L132").Value = 10621.5
$ws.Range("M132").Value = -6943.6001
$ws.Range("N132").Value = -15681.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 33555
$ws.Range("I88").Value = 33555
$ws.Range("K88").Value = 33555
$ws.Range("M88").Value = -33127

$ws.Range("H91").Value = 33555
$ws.Range("I91").Value = 33555
$ws.Range("K91").Value = 33555
$ws.Range("M91").Value = -32073

$ws.Range("H92").Value = 71333.336
$ws.Range("J92").Value = 71333.336
$ws.Range("L92").Value = 71333.336
$ws.Range("N92").Value = -76325.336

$ws.Range("H109").Value = 10000
$ws.Range("J109").Value = 10000
$ws.Range("L109").Value = 10000
$ws.Range("N109").Value = -12774

$ws.Range("H122").Value = 9424
$ws.Range("J122").Value = 10056.714
$ws.Range("L122").Value = 30170.142
$ws.Range("N122").Value = -35070.142

$ws.Range("H136").Value = 30232.459
$ws.Range("I136").Value = 40575.29
$ws.Range("K136").Value = 121725.87
$ws.Range("M136").Value = -119175.87

$ws.Range("H139").Value = 55571
$ws.Range("I139").Value = 44500
$ws.Range("J139").Value = 59999.4
$ws.Range("K139").Value = 44500
$ws.Range("L139").Value = 59999.4
$ws.Range("M139").Value = -39360
$ws.Range("N139").Value = -70279.39999999999

$ws.Range("H140").Value = 87340.28999999999
$ws.Range("J140").Value = 96658.39999999999
$ws.Range("L140").Value = 96658.39999999999
$ws.Range("N140").Value = -107018.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 6133.3335
$ws.Range("I58").Value = 6133.3335
$ws.Range("K58").Value = 6133.3335
$ws.Range("M58").Value = -5825.3335

$ws.Range("H113").Value = 1619
$ws.Range("I113").Value = 958.06665
$ws.Range("J113").Value = 2202.1765
$ws.Range("K113").Value = 2874.19995
$ws.Range("L113").Value = 6606.529500000001
$ws.Range("M113").Value = -704.1999500000002
$ws.Range("N113").Value = -10946.5295

$ws.Range("H132").Value = 24073358
$ws.Range("I132").Value = 32263614
$ws.Range("J132").Value = 991727.8
$ws.Range("K132").Value = 96790842
$ws.Range("L132").Value = 2975183.4
$ws.Range("M132").Value = -96788312
$ws.Range("N132").Value = -2980243.4
